# Apply the "Artisan keyboard/curve/events commands" documentation update
# to the Commands sheet of eventbuttons.xlsx.
#
# Five new rows are inserted right before the existing "RC Command" section
# (which starts at row 99), describing four new Artisan commands:
#   keyboard(<bool>), showCurve(<name>,<bool>), showExtraCurve(...),
#   showEvents(<event_type>,<bool>) and showBackgroundEvents(<bool>).
# Everything that used to live at row 99 onwards (RC Command, WebSocket
# Command, ...) simply shifts down by 5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Commands" sheet
$ws.Activate()

# Make room for the 5 new rows right above the old row 99 ("RC Command").
$ws.Range("A99:A103").EntireRow.Insert()

# Row 99: keyboard(<bool>)
$ws.Cells.Item(99, 2).Value2 = "keyboard(<bool>)"
$ws.Cells.Item(99, 3).Value2 = "enables/disables keyboard mode"

# Row 100: showCurve(<name>,<bool>)
$ws.Cells.Item(100, 2).Value2 = "showCurve(<name>,<bool>)"
$ws.Cells.Item(100, 3).Value2 = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

# Row 101: showExtraCurve(<extra_device>,<curve>,<bool>)
$ws.Cells.Item(101, 2).Value2 = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Cells.Item(101, 3).Value2 = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

# Row 102: showEvents(<event_type>, <bool>)
$ws.Cells.Item(102, 2).Value2 = "showEvents(<event_type>, <bool>)"
$ws.Cells.Item(102, 3).Value2 = "shows/hides the events of <event_type> in [1,..,5]"

# Row 103: showBackgroundEvents(<bool>)
$ws.Cells.Item(103, 2).Value2 = "showBackgroundEvents(<bool>)"
$ws.Cells.Item(103, 3).Value2 = "shows/hides the events of the background profile"

# Match the compact auto-fit row height used by the other two-column
# "continuation" rows in this section.
$ws.Range("A99:A103").EntireRow.RowHeight = 13.8

# Column C got a little wider to fit the new descriptions.
$ws.Columns.Item(3).ColumnWidth = 27.02

# Leave the sheet scrolled/selected roughly where the author left it.
$ws.Range("C101").Select()
